$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.726.55"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.628.64"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "595.89"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").Value = "150.27"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "0.109"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").Value = "5.69"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("D12").Value = "0.151"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").Value = "27.76"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "3.099.08"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "63.544.48"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "0.0000150"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "2.629.54"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "12.34"
$ws.Range("E18").Value = "  +7.06%  "
$ws.Range("D19").Value = "4.66"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").Value = "347.86"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "6.88"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "5.71"
$ws.Range("E23").Value = "  +2.27%  "
$ws.Range("D24").Value = "66.40"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("E25").Value = "  +11.90%  "
$ws.Range("D26").Value = "1.69"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").Value = "9.19"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("D28").Value = "572.99"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").Value = "8.21"
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "2.05"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").Value = "0.0₃0845"
$ws.Range("E33").Value = "  +2.99%  "
$ws.Range("D34").Value = "1.75"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").Value = "5.27"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").Value = "168.90"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "0.409"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "1.95"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "19.38"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("D42").Value = "168.74"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").Value = "39.88"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "3.93"
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("D45").Value = "0.0602"
$ws.Range("E45").Value = "  +5.20%  "
$ws.Range("D46").Value = "21.40"
$ws.Range("E46").Value = "  -3.87%  "
$ws.Range("D47").Value = "0.629"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "0.0250"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").Value = "1.99"
$ws.Range("E49").Value = "  +5.20%  "
$ws.Range("D50").Value = "0.0966"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "19.21"
$ws.Range("E51").Value = "  +1.63%  "
